# Updated queries for C3DC first half testcases.
#
# The workbook contains several SQL queries (stored in cells C2, B2, B3, B4,
# B5, B6, B7 of Sheet1) that joined df_participant/df_diagnoses/etc using the
# generic "id" column. They are updated here to join on the more specific
# "study_id" / "participant_id" columns instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells on Sheet1 that hold one of the SQL queries needing the JOIN fix.
$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellAddr in $cells) {
    $cell = $ws.Range($cellAddr)
    $text = $cell.Value2

    if ($text -eq $null) { continue }

    $text = $text.Replace(
        'df_participant prt ON std.id = prt."study.id"',
        'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace(
        'df_diagnoses dgn ON prt.id = dgn."participant.id"',
        'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace(
        'df_treatments trt ON prt.id = trt."participant.id"',
        'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace(
        'df_treatment_resp trr ON prt.id = trr."participant.id"',
        'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace(
        'df_survival srv ON prt.id = srv."participant.id"',
        'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace(
        'df_reference_files rfs ON std.id = rfs."study.id"',
        'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $cell.Value2 = $text
}

# Update the saved view state: the sheet had scrolled down (top-left visible
# cell becomes A6) and the active cell / selection moves from B2 to C7.
$ws.Range("C7").Select()
$win = $excel.ActiveWindow
if ($win -ne $null) {
    $win.ScrollRow = 6
    $win.ScrollColumn = 1
}
